# Moves the "System" entry to the front of the comma-separated "Recorded By"
# list in column G, preserving the relative order of the remaining entries.
# Cells where "System" is not present, or is already the first entry, are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $text = $val.ToString()
    if (-not $text.Contains("System")) {
        continue
    }

    $parts = $text -split ", "

    $systemIndex = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i].Equals("System")) {
            $systemIndex = $i
        }
    }

    if ($systemIndex -gt 0) {
        $newParts = New-Object System.Collections.ArrayList
        $newParts.Add("System") | Out-Null
        for ($i = 0; $i -lt $parts.Count; $i++) {
            if ($i -ne $systemIndex) {
                $newParts.Add($parts[$i]) | Out-Null
            }
        }
        $cell.Value2 = ($newParts -join ", ")
    }
}
